$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Current as of:" date in B1
$ws.Range("B1").Value = 44588

# Row 12 (task 7.02): now Completed, 100%
$ws.Range("B12").Value = "Completed"
$ws.Range("D12").Value = 1

# Row 13 (task 7.03): now In Progress, 0% complete
$ws.Range("B13").Value = "In Progress"
$ws.Range("D13").NumberFormat = $ws.Range("D12").NumberFormat
$ws.Range("D13").Value = 0
